$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Update the "Status" column (C) for the existing bug rows: mark issues as resolved,
# note one as "confused ??", and add a reviewer comment in a brand-new row 8.
$ws.Range("C2").Value = "solved"
$ws.Range("C4").Value = "solved"
$ws.Range("C7").Value = "confused ??"
$ws.Range("C8").Value = "on creating new roles we have to choose group so I think choose group is ok"
$ws.Range("C9").Value = "solved "
$ws.Range("C11").Value = "solved"

# Move the active selection down one row, to reflect the newly added row 8.
$ws.Range("D12").Select() | Out-Null
